$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = "'22.460.47"
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.Value = "'  +9.12%  "
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.Value = "'1.604.09"
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.Value = "'  +8.37%  "
$cell.Style = "Normal"
$cell = $ws.Range("E4")
$cell.Value = "'  -0.75%  "
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.Value = "'304.18"
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.Value = "'  +8.65%  "
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.Value = "'0.9910"
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.Value = "'  +2.03%  "
$cell.Style = "Normal"
$cell = $ws.Range("D7")
$cell.Value = "'0.3699"
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.Value = "'  +1.07%  "
$cell.Style = "Normal"
$cell = $ws.Range("D8")
$cell.Value = "'0.3393"
$cell.Style = "Normal"
$cell = $ws.Range("E8")
$cell.Value = "'  +10.11%  "
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.Value = "'42.33"
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.Value = "'1.140"
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.Value = "'  +7.51%  "
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.Value = "'0.07067"
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.Value = "'  +6.07%  "
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.Value = "'19.75"
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.Value = "'  +9.12%  "
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.Value = "'5.939"
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.Value = "'  +7.49%  "
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.Value = "'6.642"
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.Value = "'  +6.94%  "
$cell.Style = "Normal"
$cell = $ws.Range("D16")
$cell.Value = "'0.00001089"
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.Value = "'  +5.83%  "
$cell.Style = "Normal"
$cell = $ws.Range("D17")
$cell.Value = "'1.599.87"
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.Value = "'  +7.84%  "
$cell.Style = "Normal"
$cell = $ws.Range("D18")
$cell.Value = "'0.9912"
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.Value = "'  +2.02%  "
$cell.Style = "Normal"
$cell = $ws.Range("D19")
$cell.Value = "'0.06805"
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.Value = "'  +14.25%  "
$cell.Style = "Normal"
$cell = $ws.Range("D20")
$cell.Value = "'77.99"
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.Value = "'  +11.67%  "
$cell.Style = "Normal"
$cell = $ws.Range("D21")
$cell.Value = "'16.12"
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.Value = "'  +11.04%  "
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.Value = "'6.035"
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.Value = "'  +9.82%  "
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.Value = "'11.85"
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.Value = "'  +7.11%  "
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.Value = "'22.477.66"
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.Value = "'  +8.99%  "
$cell.Style = "Normal"
$cell = $ws.Range("D25")
$cell.Value = "'2.397"
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.Value = "'  +5.65%  "
$cell.Style = "Normal"
$cell = $ws.Range("D26")
$cell.Value = "'2.544"
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.Value = "'  +20.02%  "
$cell.Style = "Normal"
$cell = $ws.Range("D27")
$cell.Value = "'150.81"
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.Value = "'  +5.95%  "
$cell.Style = "Normal"
$cell = $ws.Range("D28")
$cell.Value = "'19.61"
$cell.Style = "Normal"
$cell = $ws.Range("E28")
$cell.Value = "'  +13.48%  "
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.Value = "'1.782.29"
$cell.Style = "Normal"
$cell = $ws.Range("E29")
$cell.Value = "'  +8.30%  "
$cell.Style = "Normal"
$cell = $ws.Range("D30")
$cell.Value = "'120.87"
$cell.Style = "Normal"
$cell = $ws.Range("E30")
$cell.Value = "'  +5.92%  "
$cell.Style = "Normal"
$cell = $ws.Range("D31")
$cell.Value = "'4.159"
$cell.Style = "Normal"
$cell = $ws.Range("E31")
$cell.Value = "'  +4.52%  "
$cell.Style = "Normal"
$cell = $ws.Range("D32")
$cell.Value = "'6.110"
$cell.Style = "Normal"
$cell = $ws.Range("E32")
$cell.Value = "'  +21.34%  "
$cell.Style = "Normal"
$cell = $ws.Range("D33")
$cell.Value = "'0.9535"
$cell.Style = "Normal"
$cell = $ws.Range("E33")
$cell.Value = "'  +16.07%  "
$cell.Style = "Normal"
$cell = $ws.Range("D34")
$cell.Value = "'0.08291"
$cell.Style = "Normal"
$cell = $ws.Range("E34")
$cell.Value = "'  +3.56%  "
$cell.Style = "Normal"
$cell = $ws.Range("D35")
$cell.Value = "'1.635"
$cell.Style = "Normal"
$cell = $ws.Range("E35")
$cell.Value = "'  +5.89%  "
$cell.Style = "Normal"
$cell = $ws.Range("D36")
$cell.Value = "'5.278"
$cell.Style = "Normal"
$cell = $ws.Range("E37")
$cell.Value = "'  +15.18%  "
$cell.Style = "Normal"
$cell = $ws.Range("D38")
$cell.Value = "'1.270"
$cell.Style = "Normal"
$cell = $ws.Range("E38")
$cell.Value = "'  +3.14%  "
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.Value = "'8.604"
$cell.Style = "Normal"
$cell = $ws.Range("E39")
$cell.Value = "'  +12.56%  "
$cell.Style = "Normal"
$cell = $ws.Range("D40")
$cell.Value = "'0.06092"
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.Value = "'  +5.51%  "
$cell.Style = "Normal"
$cell = $ws.Range("D41")
$cell.Value = "'0.02228"
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.Value = "'  +8.85%  "
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.Value = "'0.2029"
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.Value = "'  +7.83%  "
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.Value = "'0.9910"
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.Value = "'  +2.04%  "
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.Value = "'0.5928"
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.Value = "'  +11.69%  "
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.Value = "'3.841"
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.Value = "'  +8.51%  "
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.Value = "'13.12"
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.Value = "'  +6.85%  "
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.Value = "'0.5710"
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.Value = "'  +9.78%  "
$cell.Style = "Normal"
$cell = $ws.Range("D48")
$cell.Value = "'127.39"
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.Value = "'  +7.24%  "
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.Value = "'1.980"
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.Value = "'  +8.55%  "
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.Value = "'0.06818"
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.Value = "'  +4.81%  "
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.Value = "'73.90"
$cell.Style = "Normal"
$cell = $ws.Range("E51")
$cell.Value = "'  +9.16%  "
$cell.Style = "Normal"
